# Renamed few transcripts. Updated the DataSheet
# Change the Speaker column (D) value "Davis" to "T" for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(9, 11, 14, 16, 19, 20, 22, 27, 28, 30, 31, 33)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 4)  # Column D = Speaker
    $cell.Value = "T"
}
